$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$replacements = @(
    @{Row=1;  Col=1; Text="67×37="},
    @{Row=1;  Col=2; Text="25×78="},
    @{Row=1;  Col=3; Text="31×45="},
    @{Row=1;  Col=4; Text="23×32="},
    @{Row=1;  Col=5; Text="24×71="},
    @{Row=5;  Col=1; Text="97×17="},
    @{Row=5;  Col=2; Text="33×36="},
    @{Row=5;  Col=3; Text="32×71="},
    @{Row=5;  Col=4; Text="21×40="},
    @{Row=5;  Col=5; Text="51×86="},
    @{Row=10; Col=1; Text="88×53="},
    @{Row=10; Col=2; Text="34×20="},
    @{Row=10; Col=3; Text="15×67="},
    @{Row=10; Col=4; Text="67×49="},
    @{Row=10; Col=5; Text="47×16="},
    @{Row=15; Col=1; Text="79×98="},
    @{Row=15; Col=2; Text="82×19="},
    @{Row=15; Col=3; Text="39×62="},
    @{Row=15; Col=4; Text="17×31="},
    @{Row=15; Col=5; Text="64×43="},
    @{Row=20; Col=1; Text="69×73="},
    @{Row=20; Col=2; Text="35×14="},
    @{Row=20; Col=3; Text="20×53="},
    @{Row=20; Col=4; Text="39×34="},
    @{Row=20; Col=5; Text="20×17="}
)

foreach ($rep in $replacements) {
    $cell = $t.Cell($rep.Row, $rep.Col)
    $cell.Range.Text = $rep.Text
}
